# Rename the original sheet and create a copy for the new data (After 98 Games)
$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "After 76 Games"

# Copy the sheet right after itself -- this becomes the new "After 98 Games" tab,
# inherits all formatting/styles, and picks up the active/selected tab automatically.
$ws1.Copy([System.Type]::Missing, $ws1)
$ws2 = $wb.Worksheets.Item(2)
$ws2.Name = "After 98 Games"

# Update the cell values on the new sheet to reflect the post-Game-98 grid.
$ws2.Range("B2").Value = 20
$ws2.Range("B3").Value = 24
$ws2.Range("B4").Value = 24
$ws2.Range("D4").ClearContents()
$ws2.Range("E5").Value = 24
$ws2.Range("B6").Value = 1
$ws2.Range("E6").ClearContents()
$ws2.Range("C7").Value = 15
$ws2.Range("E8").Value = 5
$ws2.Range("F9").Value = 36
$ws2.Range("D11").Value = 2
$ws2.Range("C14").Value = 1
$ws2.Range("E15").Value = 28
$ws2.Range("I16").Value = 50
$ws2.Range("B17").Value = 2
$ws2.Range("E17").Value = 2
$ws2.Range("I17").Value = 2
$ws2.Range("D18").Value = 19
$ws2.Range("I18").ClearContents()
$ws2.Range("F19").ClearContents()
$ws2.Range("I19").ClearContents()
$ws2.Range("C24").Value = 8
$ws2.Range("F24").Value = 2
$ws2.Range("B25").ClearContents()
$ws2.Range("I25").Value = 2
$ws2.Range("J25").Value = 16
$ws2.Range("C26").Value = 2
$ws2.Range("F26").Value = 4
$ws2.Range("F27").Value = 11
$ws2.Range("J27").Value = 3
$ws2.Range("G28").Value = 4
$ws2.Range("D30").Value = 1
$ws2.Range("F31").Value = 12
$ws2.Range("G31").Value = 1
$ws2.Range("C32").Value = 8
$ws2.Range("D32").ClearContents()
$ws2.Range("F32").Value = 3
$ws2.Range("G32").Value = 1
$ws2.Range("D33").Value = 5
$ws2.Range("F33").Value = 2
$ws2.Range("C34").Value = 19
$ws2.Range("F34").Value = 15
$ws2.Range("G34").Value = 10
$ws2.Range("H34").Value = 21
$ws2.Range("I34").Value = 5
$ws2.Range("K34").Value = 4
$ws2.Range("L34").Value = 9

# Restore the new sheet's selection to match the saved workbook state.
$ws2.Range("N7").Select()
